# Generate Report for Handoff
# Updates the localization-status workbook to reflect a new handoff:
#   - Status goes from "In Translation" to "Ready for handoff"
#   - The handoff-related timestamps are refreshed
#   - The "Status" / language-status columns are widened to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 06:57:11"

# Widen the zh-cn / de-de status columns to fit "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 06:57:03"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 06:57:11"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
